# Applies the text corrections described in the diff (v1.3 -> v1.4):
# - "das Planos ... cadastradas/cadastradas" -> "dos Planos ... cadastrados" (gender/number agreement)
# - "um Planos de Capacitacao de TI" -> "um Plano de Capacitacao de TI" (singular)
# - "com a Capacitacao de TI excluida" -> "com a Capacitacao de TI nao excluida" (logic fix for the cancel-deletion case)
# - verb conjugation fixes: "escolha"->"escolhe", "selecione"->"seleciona", "preencha"->"preenche"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldA = "SYSTEM exibe a listagem das Planos de Capacitacao de TI cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$newA = "SYSTEM exibe a listagem dos Planos de Capacitacao de TI cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$cellsA = @("D10", "D20", "D30", "D40", "D57", "D72", "D87", "D103")
foreach ($addr in $cellsA) {
    $ws.Range($addr).Value = $newA
}

$oldB = "Lider de Pessoas seleciona um Planos de Capacitacao de TI da listagem"
$newB = "Lider de Pessoas seleciona um Plano de Capacitacao de TI da listagem"
$cellsB = @("B11", "B21", "B31", "B41", "B104")
foreach ($addr in $cellsB) {
    $ws.Range($addr).Value = $newB
}

$oldC = "SYSTEM exibe a listagem das Planos de Capacitacao de TI com a Capacitacao de TI excluida"
$newC = "SYSTEM exibe a listagem dos Planos de Capacitacao de TI com a Capacitacao de TI nao excluida"
$ws.Range("D13").Value = $newC

$oldD = "SYSTEM exibe a listagem das Planos de Capacitacao de TI sem a Capacitacao de TI excluida"
$newD = "SYSTEM exibe a listagem dos Planos de Capacitacao de TI sem a Capacitacao de TI excluida"
$ws.Range("D23").Value = $newD

$oldE = "SYSTEM exibe a listagem das Planos de Capacitacao de TI cadastradas apenas para visualizacao com a opcao 'Ajuda'"
$newE = "SYSTEM exibe a listagem dos Planos de Capacitacao de TI cadastrados apenas para visualizacao com a opcao 'Ajuda'"
$ws.Range("D50").Value = $newE

$oldF = "Lider de Pessoas escolha o 'Periodo Avaliativo' apropriado no campo de selecao"
$newF = "Lider de Pessoas escolhe o 'Periodo Avaliativo' apropriado no campo de selecao"
$cellsF = @("B59", "B74", "B89", "B106")
foreach ($addr in $cellsF) {
    $ws.Range($addr).Value = $newF
}

$oldG = "Lider de Pessoas selecione a 'Unidade' correspondente no campo de selecao de unidade"
$newG = "Lider de Pessoas seleciona a 'Unidade' correspondente no campo de selecao de unidade"
$cellsG = @("B60", "B75", "B90", "B107")
foreach ($addr in $cellsG) {
    $ws.Range($addr).Value = $newG
}

$oldH = "Lider de Pessoas preencha o campo 'Possiveis Capacitacoes' com informacoes sobre capacitacoes adicionais"
$newH = "Lider de Pessoas preenche o campo 'Possiveis Capacitacoes' com informacoes sobre capacitacoes adicionais"
$cellsH = @("B61", "B76", "B91", "B108")
foreach ($addr in $cellsH) {
    $ws.Range($addr).Value = $newH
}

$oldI = "Lider de Pessoas preencha o campo 'Observacao' com informacoes adicionais ou relevantes sobre o plano de capacitacao"
$newI = "Lider de Pessoas preenche o campo 'Observacao' com informacoes adicionais ou relevantes sobre o plano de capacitacao"
$cellsI = @("B62", "B77", "B92", "B109")
foreach ($addr in $cellsI) {
    $ws.Range($addr).Value = $newI
}
